$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# First, introduce every brand-new shared string in exactly the order
# they must appear in the rebuilt shared-strings table (indices 8-14):
#   8 StartTime, 9 2D, 10 Immediately, 11 5W,
#   12 0x410da55D..., 13 0xFc3BD8d2..., 14 0xEDAAb775...
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "StartTime"
$ws.Range("C3").Value = "2D"
$ws.Range("B6").Value = "Immediately"
$ws.Range("C6").Value = "5W"
$ws.Range("A4").Value = "0x410da55D45bE4f9e0616F2a1Cac7917d1baB92e5"
$ws.Range("A5").Value = "0xFc3BD8d2F74262B4dc52904cDce974822AA812bF"
$ws.Range("A6").Value = "0xEDAAb775b37A5f5098390A1e5bA2e3f6B423AE7b"

# ---------------------------------------------------------------------
# Now fill in the remaining cells (existing strings get reused slots,
# numbers are not shared strings at all).
# ---------------------------------------------------------------------
# Row 1 (header) - A1 stays "Address" (untouched)
$ws.Range("C1").Value = "Period"
$ws.Range("D1").Value = "Amount"

# Row 2 - A2 stays "0x2B0ab279A120E81C8731FEADf85262312A80897c" (untouched)
$ws.Range("B2").Value = 44214.333333333336
$ws.Range("C2").Value = "1M"
$ws.Range("D2").Value = 500000

# Row 3 - A3 stays "0xAA1D9000dDdC91227ef748C9389cB90bc9D355cF" (untouched)
$ws.Range("B3").Value = 44215.723611111112
$ws.Range("D3").Value = 10000

# Row 4
$ws.Range("B4").Value = 44217.333333333336
$ws.Range("C4").Value = "5H"
$ws.Range("D4").Value = 20000

# Row 5
$ws.Range("B5").Value = 44257.305555555555
$ws.Range("C5").Value = "3M"
$ws.Range("D5").Value = 100000

# Row 6
$ws.Range("D6").Value = 70000

# Row 7 (new, empty but styled)
$ws.Range("A7").Value = ""

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 21.5
$ws.Columns.Item(3).ColumnWidth = 12

# --- Styling: center alignment for the whole used range (row by row so we
#     don't accidentally materialise cells beyond what's populated) ---
$ws.Range("A1:D6").HorizontalAlignment = -4108
$ws.Range("A7").HorizontalAlignment = -4108

# --- Number format (date/time) for column B rows 2-6 ---
$ws.Range("B2:B6").NumberFormat = "m/d/yy h:mm"

# --- Selection / active cell ---
$ws.Range("D7").Select()
